$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.556726865544064
$ws.Range("E2").Value = 2.200360373007585
$ws.Range("F2").Value = 2.842289131771047
$ws.Range("G2").Value = 3.445728207270303
$ws.Range("H2").Value = 3.991566014111623
$ws.Range("I2").Value = 4.471761563912005
$ws.Range("J2").Value = 4.884673583290253
$ws.Range("K2").Value = 5.231691300239384
$ws.Range("L2").Value = 5.515143212015332
$ws.Range("M2").Value = 5.729462849754614
$ws.Range("N2").Value = 5.877708722108288
$ws.Range("O2").Value = 5.960482037637803
$ws.Range("P2").Value = 5.975564345232804
$ws.Range("Q2").Value = 5.935066795923256
$ws.Range("R2").Value = 5.863205667689932
$ws.Range("S2").Value = 5.775882371081892
$ws.Range("T2").Value = 5.683334764219688
$ws.Range("U2").Value = 5.591954952697852
$ws.Range("V2").Value = 5.505538978276228
$ws.Range("W2").Value = 5.426147539579697
$ws.Range("X2").Value = 5.35469906811995
$ws.Range("Y2").Value = 5.291377729038081
$ws.Range("Z2").Value = 5.235912890703711
$ws.Range("AA2").Value = 5.187769056446954
$ws.Range("AB2").Value = 5.146273340010618
$ws.Range("AC2").Value = 5.110699415919458
$ws.Range("AD2").Value = 5.080321247937529
$ws.Range("AE2").Value = 5.054445976265581
$ws.Range("AF2").Value = 5.036283868053141

$wb.Save()
